$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172; existing rows 172:193 shift down to 173:194
$ws.Rows.Item(172).Insert()

# Populate the new row 172 with the new record
$ws.Range("A172").Value = 9
$ws.Range("B172").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 44474
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 100112052
$ws.Range("G172").Value = "Albahaca"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 52
$ws.Range("K172").Value = 7000
$ws.Range("L172").Value = 7000
$ws.Range("M172").Value = 7000
$ws.Range("N172").Value = "$/docena de matas"
$ws.Range("O172").Value = "Provincia de Chacabuco"
$ws.Range("P172").Value = 1167
$ws.Range("Q172").Value = 6
$ws.Range("R172").Value = "Hortaliza"
